$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column E: refGenome ---
$ws.Range("E1").Value = "refGenome"

# hg19 for the existing hic rows (2-7)
$ws.Range("E2").Value = "hg19"
$ws.Range("E3").Value = "hg19"
$ws.Range("E4").Value = "hg19"
$ws.Range("E5").Value = "hg19"
$ws.Range("E6").Value = "hg19"
$ws.Range("E7").Value = "hg19"

# --- New rows 8-10: additional 4DN hic datasets ---
$ws.Range("C8").Value = "HFFc6"
$ws.Range("C9").Value = "H1-hESC"
$ws.Range("C10").Value = "HFF-hTERT"

$ws.Range("E8").Value = "hg38"
$ws.Range("E9").Value = "hg38"
$ws.Range("E10").Value = "hg38"

$ws.Range("A8").Value = "https://data.4dnucleome.org/files-processed/4DNFIFLJLIS5/@@download/4DNFIFLJLIS5.hic"
$ws.Range("A9").Value = "https://data.4dnucleome.org/files-processed/4DNFIOX3BGNE/@@download/4DNFIOX3BGNE.hic"
$ws.Range("A10").Value = "https://data.4dnucleome.org/files-processed/4DNFIZ4F74QR/@@download/4DNFIZ4F74QR.hic"

$ws.Range("B8").Value = "hic"
$ws.Range("B9").Value = "hic"
$ws.Range("B10").Value = "hic"

$ws.Range("D8").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("D10").Value = 0

# --- Update the selection to match the saved workbook state ---
$ws.Range("A15").Select() | Out-Null
